# Updated cryptos list on Mon Nov 27 20:53:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.876.90"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "1.995.70"
$ws.Range("E3").Value = "  -3.99%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'223.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.78%  "
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'54.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.22%  "
$ws.Range("D9").Value = "'0.375"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  -4.68%  "
$ws.Range("D12").Value = "2.289.09"
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("D13").Value = "'13.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.55%  "
$ws.Range("E14").Value = "  -6.83%  "
$ws.Range("D15").Value = "'0.732"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.64%  "
$ws.Range("E16").Value = "  -6.00%  "
$ws.Range("D17").Value = "1.985.94"
$ws.Range("E17").Value = "  -4.43%  "
$ws.Range("D18").Value = "36.862.17"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "'6.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "'68.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").Value = "0.0₃0806"
$ws.Range("E21").Value = "  -2.75%  "
$ws.Range("D22").Value = "'221.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").Value = "'2.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.63%  "
$ws.Range("D26").Value = "'164.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("E27").Value = "  -8.53%  "
$ws.Range("E28").Value = "  -4.06%  "
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("E30").Value = "  -7.76%  "
$ws.Range("E31").Value = "  -4.32%  "
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("E33").Value = "  -3.80%  "
$ws.Range("E34").Value = "  -5.93%  "
$ws.Range("E35").Value = "  -9.24%  "
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").Value = "'3.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.21%  "
$ws.Range("D39").Value = "'5.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "1.460.05"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("E41").Value = "  -6.11%  "
$ws.Range("D42").Value = "'94.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.25%  "
$ws.Range("E43").Value = "  -5.43%  "
$ws.Range("D44").Value = "'16.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("D45").Value = "'2.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.25%  "
$ws.Range("E46").Value = "  -7.88%  "
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("D50").Value = "2.178.44"
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("E51").Value = "  -4.68%  "
